# Update Calr-Itgav NATMI sheet with new TPM-derived values for the "ECs" cluster.
# The sheet lists, for every Sending-cluster x Target-cluster pair, ligand (Calr)
# expression stats for the sending cluster and receptor (Itgav) expression stats
# for the target cluster, plus several values that are derived from those two
# numbers (specificities and edge weights). Only the raw "ECs" ligand/receptor
# expression numbers changed (re-computed from the new TPM matrix); everything
# else in the sheet is recomputed from those numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand (Calr) average / total expression for cells sent from the ECs cluster.
$newLigandAvgECs   = 39.41161066666667
$newLigandTotalECs = 118.234832

# New receptor (Itgav) average / total expression for cells in the ECs cluster (target).
$newReceptorAvgECs   = 3.759736666666667
$newReceptorTotalECs = 11.27921

# Ligand average/total expression per sending cluster (column G / H), keyed by row.
# Rows 2-4 are sent from ECs, 5-7 from FAPs, 8-10 from MuSCs.
$ligandAvg = @{
    2 = $newLigandAvgECs;  3 = $newLigandAvgECs;  4 = $newLigandAvgECs
    5 = 124.0161413333333; 6 = 124.0161413333333; 7 = 124.0161413333333
    8 = 31.06188766666667; 9 = 31.06188766666667; 10 = 31.06188766666667
}
$ligandTotal = @{
    2 = $newLigandTotalECs; 3 = $newLigandTotalECs; 4 = $newLigandTotalECs
    5 = 372.048424;         6 = 372.048424;         7 = 372.048424
    8 = 93.18566300000001;  9 = 93.18566300000001;  10 = 93.18566300000001
}

# Receptor average/total expression per target cluster (column M / N), keyed by row.
# Target order within each sending-cluster block is ECs, FAPs, MuSCs.
$receptorAvg = @{
    2 = $newReceptorAvgECs; 3 = 35.81943766666667; 4 = 15.40769666666667
    5 = $newReceptorAvgECs; 6 = 35.81943766666667; 7 = 15.40769666666667
    8 = $newReceptorAvgECs; 9 = 35.81943766666667; 10 = 15.40769666666667
}
$receptorTotal = @{
    2 = $newReceptorTotalECs; 3 = 107.458313; 4 = 46.22309
    5 = $newReceptorTotalECs; 6 = 107.458313; 7 = 46.22309
    8 = $newReceptorTotalECs; 9 = 107.458313; 10 = 46.22309
}

$rows = 2..10

# Write the raw ligand/receptor average & total expression values.
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value  = $ligandAvg[$r]     # column G
    $ws.Cells.Item($r, 8).Value  = $ligandTotal[$r]   # column H
    $ws.Cells.Item($r, 13).Value = $receptorAvg[$r]   # column M
    $ws.Cells.Item($r, 14).Value = $receptorTotal[$r] # column N
}

# Ligand derived specificity = ligand value for this row's sending cluster divided by
# the sum of ligand values across the (unique) sending clusters represented.
$sendingRows = @{2 = $true; 5 = $true; 8 = $true}
$ligandAvgSum = 0.0
$ligandTotalSum = 0.0
foreach ($r in $sendingRows.Keys) { $ligandAvgSum += $ligandAvg[$r]; $ligandTotalSum += $ligandTotal[$r] }

$targetRows = @{2 = $true; 3 = $true; 4 = $true}
$receptorAvgSum = 0.0
$receptorTotalSum = 0.0
foreach ($r in $targetRows.Keys) { $receptorAvgSum += $receptorAvg[$r]; $receptorTotalSum += $receptorTotal[$r] }

foreach ($r in $rows) {
    $ligSpecAvg   = $ligandAvg[$r]   / $ligandAvgSum
    $ligSpecTotal = $ligandTotal[$r] / $ligandTotalSum
    $recSpecAvg   = $receptorAvg[$r]   / $receptorAvgSum
    $recSpecTotal = $receptorTotal[$r] / $receptorTotalSum

    $ws.Cells.Item($r, 9).Value  = $ligSpecAvg    # column I
    $ws.Cells.Item($r, 10).Value = $ligSpecTotal  # column J
    $ws.Cells.Item($r, 15).Value = $recSpecAvg    # column O
    $ws.Cells.Item($r, 16).Value = $recSpecTotal  # column P

    # Edge weights / derived specificities combine ligand (sending) and receptor (target).
    $ws.Cells.Item($r, 17).Value = $ligandAvg[$r]   * $receptorAvg[$r]    # column Q
    $ws.Cells.Item($r, 18).Value = $ligandTotal[$r] * $receptorTotal[$r] # column R
    $ws.Cells.Item($r, 19).Value = $ligSpecAvg   * $recSpecAvg    # column S
    $ws.Cells.Item($r, 20).Value = $ligSpecTotal * $recSpecTotal # column T
}
